$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Last status check on" timestamp in F1
$ws.Range("F1").Value = "Last status check on: 16.02.2022 12:45"

# D3: change from text "+0.6" to numeric 0.6
$ws.Range("D3").Value = 0.6

# E3: change from text timestamp to numeric Excel date/time value with date style
$ws.Range("E3").Value = 44608.52216435185
$ws.Range("E3").NumberFormat = "YYYY-MM-DD HH:MM:SS"
